$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-11-02 Thursday" "2023-11-03 Friday"

Replace-Text "41×34=1394" "23×38=874"
Replace-Text "39×67=2613" "51×40=2040"
Replace-Text "48×72=3456" "91×11=1001"
Replace-Text "30×44=1320" "18×56=1008"
Replace-Text "90×41=3690" "51×42=2142"

Replace-Text "30×59=1770" "32×88=2816"
Replace-Text "64×26=1664" "51×84=4284"
Replace-Text "16×84=1344" "95×56=5320"
Replace-Text "54×67=3618" "43×57=2451"
Replace-Text "79×31=2449" "77×15=1155"

Replace-Text "25×40=1000" "17×20=340"
Replace-Text "32×98=3136" "91×62=5642"
Replace-Text "13×72=936" "56×23=1288"
Replace-Text "93×32=2976" "38×15=570"
Replace-Text "42×88=3696" "59×49=2891"

Replace-Text "14×35=490" "46×54=2484"
Replace-Text "58×43=2494" "57×85=4845"
Replace-Text "84×38=3192" "80×35=2800"
Replace-Text "48×56=2688" "94×81=7614"
Replace-Text "54×41=2214" "77×29=2233"

Replace-Text "52×89=4628" "41×19=779"
Replace-Text "12×74=888" "43×49=2107"
Replace-Text "64×67=4288" "48×12=576"
Replace-Text "77×61=4697" "40×75=3000"
Replace-Text "50×68=3400" "82×44=3608"
